$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: shrink its height (100.8pt, was 115.2pt) ---
$ws.Rows.Item(7).RowHeight = 100.8

# --- New row 13: "Set different Screen Size (Tablet Size)" feature row ---
# Copy formatting from row 5, which already has the exact style pattern
# this new row needs (centered/bordered/filled A,B,D; wrapped C,E,F).
$ws.Range("A5:F5").Copy()
$ws.Range("A13:F13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A13").Value = 5
$ws.Range("B13").Value = "Set different Screen Size (Tablet Size)"
$ws.Range("C13").Value = "If the device is Tablet, use fragment to show word and meaning`ninstead of bottom sheet dialog"
$ws.Range("D13").Value = "Yes"
$ws.Range("E13").Value = "Create new xml files, adapter, and fragement for tablet size.`nIn tablet motitors, clicking word shows Meaning and Part of Speeches`nat Fragment rightside of monitor"
$ws.Range("F13").Value = "MainActivity`nBottomSheetDialogAdapter`nFragmentRecyclerViewAdapter`nMainActivityRecyclerViewAdapter`nWordBottomSheetDialog`nMainFragment`nOnItemClick (interface)`nenter_from_right.xml (anim)`nactivity_main.xml`nactivity_main_fragment_xml`nrv_fragment_pos_mean.xml"

$ws.Rows.Item(13).RowHeight = 187.2

# --- Selection moves to B11 ---
[void]$ws.Range("B11").Select()
